# "Add files via upload" - refresh of the daily COVID-19 Valais figures.
# The source workbook got a round of corrections/updates to the tail of the
# data table (rows 168-172, 213-234) plus three brand-new rows worth of
# data that used to be blank (row 234), and the frozen-pane / active-cell
# view state moved along with the newly-entered data (row 234, col H).
#
# Columns N/J/K/B/H are driven by "= shared" formulas that recompute on
# their own once the precedent cells below are written, so only the plain
# input cells are touched here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Quarantine/contact figures corrected for late Aug/early Sep (rows 168-172) ---
$ws.Range("O168").Value = 192
$ws.Range("O169").Value = 208
$ws.Range("O170").Value = 215
$ws.Range("O171").Value = 230
$ws.Range("O172").Value = 220

# --- Isolation (N) / quarantine (O) revisions, early-mid October (rows 213-227) ---
$ws.Range("N213").Value = 107
$ws.Range("N214").Value = 93
$ws.Range("N215").Value = 92
$ws.Range("N216").Value = 78
$ws.Range("N217").Value = 95
$ws.Range("N218").Value = 101
$ws.Range("N219").Value = 128
$ws.Range("N220").Value = 132

$ws.Range("N221").Value = 134
$ws.Range("O221").Value = 416

$ws.Range("N222").Value = 142
$ws.Range("O222").Value = 482

$ws.Range("N223").Value = 172
$ws.Range("O223").Value = 464

$ws.Range("N224").Value = 225
$ws.Range("O224").Value = 525

$ws.Range("N225").Value = 291
$ws.Range("O225").Value = 656

$ws.Range("N226").Value = 384
$ws.Range("O226").Value = 749

$ws.Range("N227").Value = 503
$ws.Range("O227").Value = 789

# --- Row 228 (2020-10-10): new admission + non-ICU hospitalisation counts ---
$ws.Range("D228").Value = 3
$ws.Range("G228").Value = 18
$ws.Range("N228").Value = 552
$ws.Range("O228").Value = 820
$ws.Range("P228").Value = 897

# --- Row 229 (2020-10-11) ---
$ws.Range("C229").Value = 68
$ws.Range("G229").Value = 21
$ws.Range("N229").Value = 597
$ws.Range("O229").Value = 859
$ws.Range("P229").Value = 913

# --- Row 230 (2020-10-12) ---
$ws.Range("G230").Value = 24
$ws.Range("N230").Value = 769
$ws.Range("O230").Value = 892
$ws.Range("P230").Value = 894

# --- Row 231 (2020-10-13) ---
$ws.Range("C231").Value = 157
$ws.Range("G231").Value = 28
$ws.Range("N231").Value = 898
$ws.Range("O231").Value = 775
$ws.Range("P231").Value = 846

# --- Row 232 (2020-10-14) ---
$ws.Range("C232").Value = 95
$ws.Range("D232").Value = 2
$ws.Range("G232").Value = 29
$ws.Range("I232").Value = 1
$ws.Range("N232").Value = 958
$ws.Range("O232").Value = 598
$ws.Range("P232").Value = 739

# --- Row 233 (2020-10-15) ---
$ws.Range("C233").Value = 7
$ws.Range("G233").Value = 28
$ws.Range("L233").Value = 1
$ws.Range("N233").Value = 914
$ws.Range("O233").Value = 505
$ws.Range("P233").Value = 659

# --- Row 234 (2020-10-16): was a blank placeholder row, now filled in ---
$ws.Range("C234").Value = 0
$ws.Range("D234").Value = 0
$ws.Range("E234").Value = 3
$ws.Range("F234").Value = 1
$ws.Range("G234").Value = 28
$ws.Range("I234").Value = 0
$ws.Range("L234").Value = 0
$ws.Range("M234").Value = 0
$ws.Range("N234").Value = 868
$ws.Range("O234").Value = 365
$ws.Range("P234").Value = 605

# --- View state: active selection followed the newly-entered data ---
$ws.Range("H234").Select()
